$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add new "up/down" analysis columns (X = change, Y = Up/Down verdict) ---
$ws.Range("X3").Value = -0.93999500000001035
$ws.Range("Y3").Value = "Down"

# --- Row 4: new trade record appended to the sentiment/trading log ---
$ws.Range("A4").Value = 42633.890532407408
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 36
$ws.Range("E4").Value = 15261
$ws.Range("F4").Value = 2660
$ws.Range("G4").Value = 64
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 91
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 18783
$ws.Range("L4").Value = 341
$ws.Range("M4").Value = 161
$ws.Range("N4").Value = 55
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = "Noun"
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.86
$ws.Range("S4").Value = 0.0262
$ws.Range("S4").NumberFormat = $ws.Range("S3").NumberFormat
$ws.Range("T4").Value = -2.66
$ws.Range("U4").Value = 15.05
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0
